# Weekly update: insert a new price record for "Ají" (Inferno variety)
# at row 20, pushing the existing rows 20-39 down to 21-40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 20 (existing rows shift down).
$ws.Rows.Item(20).Insert()

# Fill in the new row 20 with the new weekly record.
$ws.Cells.Item(20, 1).Value = 11
$ws.Cells.Item(20, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value = "Bíobío"
$ws.Cells.Item(20, 4).Value2 = 44447
$ws.Cells.Item(20, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = 100112021
$ws.Cells.Item(20, 7).Value = "Ají"
$ws.Cells.Item(20, 8).Value = "Inferno"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 40
$ws.Cells.Item(20, 11).Value = 35000
$ws.Cells.Item(20, 12).Value = 36000
$ws.Cells.Item(20, 13).Value = 35500
$ws.Cells.Item(20, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(20, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(20, 16).Value = 2958
$ws.Cells.Item(20, 17).Value = 12
$ws.Cells.Item(20, 18).Value = "Hortaliza"
